# Penalty Reward System (unfinished) - trim trailing data rows and
# adjust the last remaining row's quantity on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
# Remove rows 8-13 (everything after the week of 2023-06-11)
$ws1.Range("A8:B13").EntireRow.Delete()
# The last remaining row (7) drops its quantity from 10 to 5
$ws1.Range("B7").Value = 5

# --- Sheet "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
# Remove row 7 (the July 2023 month)
$ws2.Range("A7:B7").EntireRow.Delete()
# The last remaining row (6) drops its quantity from 20 to 5
$ws2.Range("B6").Value = 5
